$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D are plain numeric-looking text (e.g. "0.0636") in the
# source data; Excel's COM type-inference would silently coerce these to
# real numbers on a bare .Value assignment. Force text storage by setting
# NumberFormat to Text ("@") before the write, then ClearFormats() to drop
# the now-unneeded style index again (these cells carry no style in the
# original workbook).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '26.160.33'
$ws.Range('E2').Value = '  +1.53%  '
Set-TextValue $ws.Range('D3') '1.642.99'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.19%  '
Set-TextValue $ws.Range('D5') '216.69'
$ws.Range('E5').Value = '  +0.41%  '
Set-TextValue $ws.Range('D6') '0.505'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +0.54%  '
Set-TextValue $ws.Range('D9') '0.0636'
$ws.Range('E9').Value = '  +0.40%  '
Set-TextValue $ws.Range('D10') '19.74'
$ws.Range('E10').Value = '  +0.91%  '
Set-TextValue $ws.Range('D11') '0.0791'
$ws.Range('E11').Value = '  -0.26%  '
Set-TextValue $ws.Range('D12') '1.870.96'
Set-TextValue $ws.Range('D13') '4.29'
$ws.Range('E13').Value = '  +0.86%  '
Set-TextValue $ws.Range('D14') '1.655.54'
$ws.Range('E14').Value = '  +1.21%  '
Set-TextValue $ws.Range('D15') '0.545'
$ws.Range('E15').Value = '  -3.05%  '
Set-TextValue $ws.Range('D16') '0.0₃0763'
$ws.Range('E16').Value = '  +0.11%  '
Set-TextValue $ws.Range('D17') '63.21'
Set-TextValue $ws.Range('D18') '26.177.52'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  -0.53%  '
Set-TextValue $ws.Range('D21') '194.87'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  +0.72%  '
Set-TextValue $ws.Range('D23') '6.27'
$ws.Range('E23').Value = '  -0.40%  '
Set-TextValue $ws.Range('D24') '1.79'
$ws.Range('E24').Value = '  -2.67%  '
$ws.Range('E25').Value = '  -0.13%  '
Set-TextValue $ws.Range('D26') '142.57'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('E27').Value = '  +1.43%  '
Set-TextValue $ws.Range('D29') '15.59'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('E31').Value = '  +1.87%  '
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('E36').Value = '  +0.82%  '
Set-TextValue $ws.Range('D37') '1.131.26'
$ws.Range('E37').Value = '  +0.10%  '
Set-TextValue $ws.Range('D38') '0.553'
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  -0.23%  '
Set-TextValue $ws.Range('D42') '100.39'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('E44').Value = '  -0.54%  '
Set-TextValue $ws.Range('D45') '1.779.74'
$ws.Range('E45').Value = '  +0.58%  '
Set-TextValue $ws.Range('D46') '0.0₆0111'
$ws.Range('E46').Value = '  -0.99%  '
Set-TextValue $ws.Range('D47') '56.81'
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('E48').Value = '  +4.31%  '
Set-TextValue $ws.Range('D49') '0.0517'
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D50') '0.417'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '7.67'
$ws.Range('E51').Value = '  +3.29%  '
